$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update precision, recall, f1-score, f2-score, NDCG rows (2-6) with numeric values
$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 0.4
$ws.Range("I2").Value = 0.25
$ws.Range("M2").Value = 1
$ws.Range("Q2").Value = 0.2
$ws.Range("S2").Value = 0.25
$ws.Range("U2").Value = 0.5
$ws.Range("W2").Value = 1
$ws.Range("C3").Value = 0.5
$ws.Range("F3").Value = 1
$ws.Range("I3").Value = 0.5
$ws.Range("M3").Value = 0.5
$ws.Range("Q3").Value = 0.5
$ws.Range("S3").Value = 0.5
$ws.Range("U3").Value = 0.5
$ws.Range("W3").Value = 0.5
$ws.Range("C4").Value = 0.6666666666666666
$ws.Range("F4").Value = 0.5714285714285715
$ws.Range("I4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.6666666666666666
$ws.Range("Q4").Value = 0.2857142857142858
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("U4").Value = 0.5
$ws.Range("W4").Value = 0.6666666666666666
$ws.Range("C5").Value = 0.5555555555555556
$ws.Range("F5").Value = 0.7692307692307692
$ws.Range("I5").Value = 0.4166666666666667
$ws.Range("M5").Value = 0.5555555555555556
$ws.Range("Q5").Value = 0.3846153846153846
$ws.Range("S5").Value = 0.4166666666666667
$ws.Range("U5").Value = 0.5
$ws.Range("W5").Value = 0.5555555555555556
$ws.Range("C6").Value = 0.8262346571285599
$ws.Range("F6").Value = 0.9639404333166532
$ws.Range("I6").Value = 0.41311732856428
$ws.Range("M6").Value = 0.8262346571285599
$ws.Range("Q6").Value = 0.52129602861432
$ws.Range("S6").Value = 0.52129602861432
$ws.Range("U6").Value = 0.8262346571285599
$ws.Range("W6").Value = 0.8262346571285599

# Update M1, M3, M5 boolean rows (7-9)
$ws.Range("C7").Value = $true
$ws.Range("F7").Value = $true
$ws.Range("M7").Value = $true
$ws.Range("U7").Value = $true
$ws.Range("W7").Value = $true
$ws.Range("C8").Value = $true
$ws.Range("F8").Value = $true
$ws.Range("I8").Value = $true
$ws.Range("M8").Value = $true
$ws.Range("Q8").Value = $true
$ws.Range("S8").Value = $true
$ws.Range("U8").Value = $true
$ws.Range("W8").Value = $true
$ws.Range("C9").Value = $true
$ws.Range("F9").Value = $true
$ws.Range("I9").Value = $true
$ws.Range("M9").Value = $true
$ws.Range("Q9").Value = $true
$ws.Range("S9").Value = $true
$ws.Range("U9").Value = $true
$ws.Range("W9").Value = $true

# Update position row (10) with numeric values
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("I10").Value = 3
$ws.Range("M10").Value = 1
$ws.Range("Q10").Value = 2
$ws.Range("S10").Value = 2
$ws.Range("U10").Value = 1
$ws.Range("W10").Value = 1
